# Update the Jogos do Dia Betfair Back/Lay sheet for 2025-10-15:
#  - Drop the Japanese J League 3 fixture (old row 2, Kanazawa vs Matsumoto)
#    and the Colombian Primera A fixture (old row 7, Deportivo Pereira vs
#    Millonarios) moves up into row 6 with refreshed odds; all Brazilian Serie A
#    fixtures shift up one row with updated odds.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 (Japanese J League 3) is dropped by deleting the current last
# row (row 7); every other row keeps its row number, and we overwrite rows 2-6
# below with each fixture's new details / odds, matching the other rows shifting
# up one position.
$ws.Rows(7).Delete()

# --- Row 2: SE Palmeiras vs Red Bull Bragantino (Brazilian Serie A) ---
$ws.Range("A2").Value = 'Brazilian Serie A'
$ws.Range("C2").Value = '19:00:00'
$ws.Range("D2").Value = 'SE Palmeiras'
$ws.Range("E2").Value = 'Red Bull Bragantino'
$ws.Range("F2").Value = 1.37
$ws.Range("G2").Value = 1.4
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 5.1
$ws.Range("K2").Value = 5.5
$ws.Range("N2").Value = 3.95
$ws.Range("O2").Value = 1.28
$ws.Range("P2").Value = 2.02
$ws.Range("Q2").Value = 1.83

# --- Row 3: Botafogo FR vs Flamengo (Brazilian Serie A) ---
$ws.Range("C3").Value = '19:30:00'
$ws.Range("D3").Value = 'Botafogo FR'
$ws.Range("E3").Value = 'Flamengo'
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 4.3
$ws.Range("H3").Value = 2.14
$ws.Range("I3").Value = 2.3
$ws.Range("J3").Value = 3.15
$ws.Range("K3").Value = 3.6
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 1.67
$ws.Range("Q3").Value = 2.28

# --- Row 4: Mirassol vs Internacional (Brazilian Serie A) ---
$ws.Range("C4").Value = '20:00:00'
$ws.Range("D4").Value = 'Mirassol'
$ws.Range("E4").Value = 'Internacional'
$ws.Range("F4").Value = 2.12
$ws.Range("G4").Value = 2.3
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 4.2
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 3.65
$ws.Range("P4").Value = 1.75
$ws.Range("Q4").Value = 1.95

# --- Row 5: Sport Recife vs Ceara SC Fortaleza (Brazilian Serie A) ---
$ws.Range("D5").Value = 'Sport Recife'
$ws.Range("E5").Value = 'Ceara SC Fortaleza'
$ws.Range("F5").Value = 2.44
$ws.Range("G5").Value = 2.68
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 3.85
$ws.Range("J5").Value = 3.05
$ws.Range("K5").Value = 3.35
$ws.Range("P5").Value = 1.58
$ws.Range("Q5").Value = 2.48

# --- Row 6: Deportivo Pereira vs Millonarios (Colombian Primera A) ---
$ws.Range("A6").Value = 'Colombian Primera A'
$ws.Range("C6").Value = '20:20:00'
$ws.Range("D6").Value = 'Deportivo Pereira'
$ws.Range("E6").Value = 'Millonarios'
$ws.Range("F6").Value = 2.46
$ws.Range("G6").Value = 2.84
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 3.85
$ws.Range("J6").Value = 3
$ws.Range("P6").Value = 1.6
$ws.Range("Q6").Value = 2.4

